# Generate Report for Handback
# Update the handoff/handback timestamps for the "a16099ae-..." row (row 4)
# on both the zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-19 04:37:15"
$wsZhCn.Range("H4").Value = "2016-03-19 04:37:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-19 04:37:18"
$wsDeDe.Range("H4").Value = "2016-03-19 04:37:42"
